# resolvido o problema da falta do rotulo titulo
# Repurpose the "Sao Joao" presentation into a "Robotica" presentation.
#
# NOTE: the host's TextRange.Text setter tries to preserve per-run
# formatting by diffing the new string against the shape's current text
# (matching a common leading/trailing substring). When the old and new
# text happen to share a character at the start/end this produces an
# extra, unwanted <a:r> split. To keep each paragraph as a single run
# (matching the canonical authoring output) we first blank out the
# text with a placeholder that shares nothing with either the old or
# new text, then assign the final text.

function Set-ShapeText($shape, [string]$text) {
    $shape.TextFrame.TextRange.Text = "#"
    $shape.TextFrame.TextRange.Text = $text
}

$p = $ppt.ActivePresentation

# Slide 1: title slide
$s1 = $p.Slides.Item(1)
Set-ShapeText $s1.Shapes.Item(1) "Robótica"
Set-ShapeText $s1.Shapes.Item(2) "Autor: Rei Robô"

# Slide 2: Introducao a Robotica
$s2 = $p.Slides.Item(2)
Set-ShapeText $s2.Shapes.Item(1) "Introdução à Robótica"
Set-ShapeText $s2.Shapes.Item(2) "Definição de Robótica`rHistória da Robótica`rAplicações da Robótica"

# Slide 3: Componentes de um Robo
$s3 = $p.Slides.Item(3)
Set-ShapeText $s3.Shapes.Item(1) "Componentes de um Robô"
Set-ShapeText $s3.Shapes.Item(2) "Sensores`rAtuadores`rControladores"

# Slide 4: Tipos de Robos
$s4 = $p.Slides.Item(4)
Set-ShapeText $s4.Shapes.Item(1) "ipos de Robôs"
Set-ShapeText $s4.Shapes.Item(2) "Robôs Industriais`rRobôs Móveis`rRobôs Autônomos`rRobôs Colaborativos"

# Slide 5: Aplicacoes da Robotica
$s5 = $p.Slides.Item(5)
Set-ShapeText $s5.Shapes.Item(1) "Aplicações da Robótica"
Set-ShapeText $s5.Shapes.Item(2) "Linha de produção automatizada`rRobôs cirúrgicos`rRobôs de exploração espacial`rRobôs de entretenimento"

# Slide 6: Robotica e Inteligencia Artificial
$s6 = $p.Slides.Item(6)
Set-ShapeText $s6.Shapes.Item(1) "Robótica e Inteligência Artificial"
Set-ShapeText $s6.Shapes.Item(2) "Integração de IA nos robôs`rAprendizado de Máquina`rRobôs autônomos"

# Slide 7: Etica na Robotica
$s7 = $p.Slides.Item(7)
Set-ShapeText $s7.Shapes.Item(1) "Ética na Robótica"
Set-ShapeText $s7.Shapes.Item(2) "Questões éticas em relação à autonomia dos robôs`rImpacto social e econômico dos robôs`rResponsabilidade e accountability na utilização de robôs"

# Slide 8: Desafios da Robotica
$s8 = $p.Slides.Item(8)
Set-ShapeText $s8.Shapes.Item(1) "Desafios da Robótica"
Set-ShapeText $s8.Shapes.Item(2) "Manipulação de objetos complexos`rNavegação autônoma em ambientes dinâmicos`rInterpretação e execução de comandos humanos"

# Slide 9: Futuro da Robotica
$s9 = $p.Slides.Item(9)
Set-ShapeText $s9.Shapes.Item(1) "Futuro da Robótica"
Set-ShapeText $s9.Shapes.Item(2) "Avanços tecnológicos`rIntegração cada vez maior de robôs na sociedade`rNovas aplicações da robótica"

# Slide 10: Conclusao (title unchanged, content merges into a single paragraph)
$s10 = $p.Slides.Item(10)
Set-ShapeText $s10.Shapes.Item(2) "A robótica é uma área em constante evolução, com aplicações em diversas áreas e um potencial de transformação da sociedade. É importante considerar os desafios éticos e sociais relacionados à utilização de robôs, buscando sempre um equilíbrio entre inovação e responsabilidade. A robótica continuará a desempenhar um papel fundamental no futuro, criando novas possibilidades e desafios para a humanidade."
